$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'300.52"
$ws.Range("E2").Value = "'-0.69%"
$ws.Range("G2").Value = "'4"
$ws.Range("D3").Value = "'38.34"
$ws.Range("E3").Value = "'9.23%"
$ws.Range("G3").Value = "'4"
$ws.Range("D4").Value = "'4.984"
$ws.Range("E4").Value = "'-3.35%"
$ws.Range("G4").Value = "'4"
$ws.Range("D5").Value = "'0.07728"
$ws.Range("E5").Value = "'-0.55%"
$ws.Range("G5").Value = "'4"
$ws.Range("D6").Value = "'2.173"
$ws.Range("E6").Value = "'-6.71%"
$ws.Range("G6").Value = "'4"
$ws.Range("D7").Value = "'7.962"
$ws.Range("E7").Value = "'-0.99%"
$ws.Range("G7").Value = "'4"
$ws.Range("D8").Value = "'3.996"
$ws.Range("E8").Value = "'1.16%"
$ws.Range("G8").Value = "'4"
$ws.Range("D9").Value = "'0.9145"
$ws.Range("E9").Value = "'-1.74%"
$ws.Range("G9").Value = "'4"
$ws.Range("D10").Value = "'0.09084"
$ws.Range("E10").Value = "'-9.51%"
$ws.Range("G10").Value = "'4"
$ws.Range("D11").Value = "'0.1798"
$ws.Range("E11").Value = "'0.55%"
$ws.Range("G11").Value = "'4"
$ws.Range("D12").Value = "'0.08491"
$ws.Range("E12").Value = "'-1.80%"
$ws.Range("G12").Value = "'4"
$ws.Range("D13").Value = "'0.03527"
$ws.Range("E13").Value = "'6.00%"
$ws.Range("G13").Value = "'4"
$ws.Range("D14").Value = "'0.09941"
$ws.Range("E14").Value = "'0.28%"
$ws.Range("G14").Value = "'4"
$ws.Range("D15").Value = "'0.001484"
$ws.Range("E15").Value = "'-0.99%"
$ws.Range("G15").Value = "'4"
$ws.Range("D16").Value = "'0.005693"
$ws.Range("E16").Value = "'-1.34%"
$ws.Range("G16").Value = "'4"
$ws.Range("D17").Value = "'3.478"
$ws.Range("E17").Value = "'0.50%"
$ws.Range("G17").Value = "'4"
$ws.Range("G18").Value = "'4"
$ws.Range("E19").Value = "'3.15%"
$ws.Range("G19").Value = "'4"
$ws.Range("E20").Value = "'-1.26%"
$ws.Range("G20").Value = "'4"
$ws.Range("D21").Value = "'4.555"
$ws.Range("E21").Value = "'6.06%"
$ws.Range("G21").Value = "'4"
$ws.Range("D22").Value = "'0.2233"
$ws.Range("E22").Value = "'-3.02%"
$ws.Range("G22").Value = "'4"
$ws.Range("D23").Value = "'0.04655"
$ws.Range("E23").Value = "'2.16%"
$ws.Range("G23").Value = "'4"
$ws.Range("D24").Value = "'0.001229"
$ws.Range("E24").Value = "'1.24%"
$ws.Range("G24").Value = "'4"
$ws.Range("D25").Value = "'0.004442"
$ws.Range("E25").Value = "'1.64%"
$ws.Range("G25").Value = "'4"
$ws.Range("D26").Value = "'0.0001301"
$ws.Range("E26").Value = "'-0.06%"
$ws.Range("G26").Value = "'4"
$ws.Range("D27").Value = "'0.0004753"
$ws.Range("E27").Value = "'40.01%"
$ws.Range("G27").Value = "'4"
$ws.Range("G28").Value = "'4"
$ws.Range("G29").Value = "'4"
$ws.Range("G30").Value = "'4"
$ws.Range("G31").Value = "'4"
$ws.Range("G32").Value = "'4"
$ws.Range("G33").Value = "'4"
$ws.Range("G34").Value = "'4"
$ws.Range("G35").Value = "'4"
$ws.Range("G36").Value = "'4"
$ws.Range("G37").Value = "'4"
$ws.Range("G38").Value = "'4"
$ws.Range("D39").Value = "'0.01744"
$ws.Range("E39").Value = "'-2.79%"
$ws.Range("G39").Value = "'4"
$ws.Range("D40").Value = "'0.04677"
$ws.Range("E40").Value = "'-2.77%"
$ws.Range("G40").Value = "'4"
$ws.Range("D41").Value = "'0.007917"
$ws.Range("E41").Value = "'1.63%"
$ws.Range("G41").Value = "'4"
$ws.Range("E42").Value = "'-1.71%"
$ws.Range("G42").Value = "'4"
$ws.Range("D43").Value = "'0.007664"
$ws.Range("E43").Value = "'11.97%"
$ws.Range("G43").Value = "'4"
$ws.Range("D44").Value = "'0.002302"
$ws.Range("E44").Value = "'9.08%"
$ws.Range("G44").Value = "'4"
$ws.Range("D45").Value = "'0.01010"
$ws.Range("E45").Value = "'6.97%"
$ws.Range("G45").Value = "'4"
$ws.Range("D46").Value = "'0.00006022"
$ws.Range("E46").Value = "'-1.49%"
$ws.Range("G46").Value = "'4"
$ws.Range("E47").Value = "'-0.05%"
$ws.Range("G47").Value = "'4"
$ws.Range("D48").Value = "'8.704"
$ws.Range("E48").Value = "'184.73%"
$ws.Range("G48").Value = "'4"
$ws.Range("E49").Value = "'34.77%"
$ws.Range("G49").Value = "'4"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("E50").Value = "'-0.05%"
$ws.Range("G50").Value = "'4"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").Value = "'-0.05%"
$ws.Range("G51").Value = "'4"
